$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.700.06"
$ws.Range("E2").Value = "  +0.89%  "

$ws.Range("D3").Value = "'3.456.37"
$ws.Range("E3").Value = "  +1.69%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'578.43"
$ws.Range("E5").Value = "  +0.93%  "

$ws.Range("D6").Value = "'145.85"
$ws.Range("E6").Value = "  +5.04%  "

$ws.Range("D7").Value = "'3.458.47"
$ws.Range("E7").Value = "  +1.76%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  +2.04%  "

$ws.Range("D10").Value = "'7.71"
$ws.Range("E10").Value = "  +0.04%  "

$ws.Range("E11").Value = "  +3.95%  "

$ws.Range("D12").Value = "'0.392"
$ws.Range("E12").Value = "  +2.90%  "

$ws.Range("D13").Value = "'4.046.86"
$ws.Range("E13").Value = "  +1.74%  "

$ws.Range("D14").Value = "'28.58"
$ws.Range("E14").Value = "  +7.21%  "

$ws.Range("E15").Value = "  -0.39%  "

$ws.Range("D17").Value = "'3.453.87"
$ws.Range("E17").Value = "  +1.67%  "

$ws.Range("D18").Value = "'61.793.95"
$ws.Range("E18").Value = "  +1.01%  "

$ws.Range("D19").Value = "'6.40"
$ws.Range("E19").Value = "  +7.64%  "

$ws.Range("D20").Value = "'14.38"
$ws.Range("E20").Value = "  +3.73%  "

$ws.Range("D21").Value = "'9.45"
$ws.Range("E21").Value = "  +1.58%  "

$ws.Range("D22").Value = "'403.93"
$ws.Range("E22").Value = "  +7.04%  "

$ws.Range("D23").Value = "'0.569"
$ws.Range("E23").Value = "  +2.93%  "

$ws.Range("D24").Value = "'74.61"
$ws.Range("E24").Value = "  +4.82%  "

$ws.Range("D25").Value = "'5.78"
$ws.Range("E25").Value = "  +0.60%  "

$ws.Range("E26").Value = "  -0.57%  "

$ws.Range("E27").Value = "  +1.82%  "

$ws.Range("D28").Value = "'3.591.60"
$ws.Range("E28").Value = "  +1.64%  "

$ws.Range("D29").Value = "'0.184"
$ws.Range("E29").Value = "  +4.77%  "

$ws.Range("D30").Value = "'7.65"
$ws.Range("E30").Value = "  +3.03%  "

$ws.Range("E31").Value = "  +0.23%  "

$ws.Range("D32").Value = "'8.28"
$ws.Range("E32").Value = "  +1.60%  "

$ws.Range("E33").Value = "  +2.15%  "

$ws.Range("E34").Value = "  -9.68%  "

$ws.Range("E35").Value = "  -0.08%  "

$ws.Range("D36").Value = "'24.01"
$ws.Range("E36").Value = "  +2.38%  "

$ws.Range("D37").Value = "'7.09"
$ws.Range("E37").Value = "  +2.66%  "

$ws.Range("D38").Value = "'3.482.41"
$ws.Range("E38").Value = "  +1.91%  "

$ws.Range("D39").Value = "'1.58"
$ws.Range("E39").Value = "  +0.37%  "

$ws.Range("D40").Value = "'5.16"
$ws.Range("E40").Value = "  +0.80%  "

$ws.Range("D41").Value = "'166.90"
$ws.Range("E41").Value = "  +0.40%  "

$ws.Range("D42").Value = "'0.0794"
$ws.Range("E42").Value = "  +2.80%  "

$ws.Range("E43").Value = "  +5.26%  "

$ws.Range("D44").Value = "'0.804"
$ws.Range("E44").Value = "  +3.40%  "

$ws.Range("E45").Value = "  +0.10%  "

$ws.Range("E46").Value = "  +2.93%  "

$ws.Range("E47").Value = "  +1.16%  "

$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  -0.03%  "

$ws.Range("D49").Value = "'2.613.67"
$ws.Range("E49").Value = "  +3.55%  "

$ws.Range("D50").Value = "'1.16"
$ws.Range("E50").Value = "  -1.16%  "

$ws.Range("D51").Value = "'6.97"
$ws.Range("E51").Value = "  +2.73%  "
